# "added jer's comments 26092013"
#
# The sheet had a merged title row (A1:D1 = "Table 1.  Primers used in this
# study for PCR amplification and sequencing reactions.") sitting above the
# header row. That title row is removed entirely, which shifts every row
# below it up by one (header row 2 -> row 1, data rows 3-10 -> rows 2-9, and
# the trailing blank-styled cell row 14 -> row 13).
#
# Separately, the primer-reference citation for the last row (16S1M) is
# corrected from "(Fu, 2000)" to "(Fu 2000)" (comma removed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole first row (the merged title banner). Excel shifts
# everything below it up by one row and drops the now-pointless merge.
$ws.Rows("1:1").Delete() | Out-Null

# Fix the reference text on the last data row (now row 9, was row 10) —
# drop the comma: "(Fu, 2000)" -> "(Fu 2000)".
$ws.Range("D9").Value = "(Fu 2000)"

# Leave the selection where the author ended up after editing.
$ws.Range("G16").Select() | Out-Null
